$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the device identifier and IP address values in row 2.
# Leading apostrophe keeps these text values tagged with the
# "quote prefix" cell style (they look numeric-ish), matching how
# the existing cells were already formatted.
$ws.Range("C2").Value = "'IND_DAU_51"
$ws.Range("E2").Value = "'10.75.58.51"

# Move the active selection to the last edited cell
$ws.Range("E2").Select()
